$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 61
$ws.Range("B61").Value = 6838387
$ws.Range("E61").Value = "PSV"
$ws.Range("F61").Value = "FC Volendam"
$ws.Range("G61").Value = 3
$ws.Range("H61").Value = 1
$ws.Range("I61").Value = 1
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = "H"
$ws.Range("L61").Value = 1.125
$ws.Range("M61").Value = 8.5
$ws.Range("N61").Value = 17
$ws.Range("O61").Value = 1.035
$ws.Range("P61").Value = 17
$ws.Range("Q61").Value = 34
$ws.Range("R61").Value = -3.75
$ws.Range("S61").Value = 2
$ws.Range("T61").Value = 1.85
$ws.Range("U61").Value = 4.75
$ws.Range("V61").Value = 1.85
$ws.Range("W61").Value = 2
$ws.Range("X61").Value = 0.03499999999999992
$ws.Range("Y61").Value = -1
$ws.Range("Z61").Value = -1
$ws.Range("AA61").Value = -1
$ws.Range("AB61").Value = 0.8500000000000001
$ws.Range("AC61").Value = -1
$ws.Range("AD61").Value = 1

# Row 62
$ws.Range("B62").Value = 6838386
$ws.Range("E62").Value = "FC Utrecht"
$ws.Range("F62").Value = "Almere City FC"
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 2
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = "A"
$ws.Range("L62").Value = 1.444
$ws.Range("M62").Value = 4.333
$ws.Range("N62").Value = 7
$ws.Range("O62").Value = 1.55
$ws.Range("P62").Value = 4
$ws.Range("Q62").Value = 6
$ws.Range("R62").Value = -1
$ws.Range("S62").Value = 1.94
$ws.Range("T62").Value = 1.96
$ws.Range("U62").Value = 3
$ws.Range("V62").Value = 2
$ws.Range("W62").Value = 1.85
$ws.Range("X62").Value = -1
$ws.Range("Y62").Value = -1
$ws.Range("Z62").Value = 5
$ws.Range("AA62").Value = -1
$ws.Range("AB62").Value = 0.96
$ws.Range("AC62").Value = -1
$ws.Range("AD62").Value = 0.8500000000000001

# Row 88
$ws.Range("B88").Value = 6838411
$ws.Range("E88").Value = "Almere City FC"
$ws.Range("F88").Value = "Go Ahead Eagles"
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = "D"
$ws.Range("L88").Value = 2.625
$ws.Range("M88").Value = 3.5
$ws.Range("N88").Value = 2.5
$ws.Range("O88").Value = 2.9
$ws.Range("P88").Value = 3.6
$ws.Range("Q88").Value = 2.3
$ws.Range("R88").Value = 0.25
$ws.Range("S88").Value = 1.875
$ws.Range("T88").Value = 1.975
$ws.Range("U88").Value = 2.75
$ws.Range("V88").Value = 1.9
$ws.Range("W88").Value = 1.95
$ws.Range("X88").Value = -1
$ws.Range("Y88").Value = 2.6
$ws.Range("Z88").Value = -1
$ws.Range("AA88").Value = 0.4375
$ws.Range("AB88").Value = -0.5
$ws.Range("AC88").Value = -1
$ws.Range("AD88").Value = 0.95

# Row 89
$ws.Range("B89").Value = 6838413
$ws.Range("E89").Value = "Sparta Rotterdam"
$ws.Range("F89").Value = "RKC"
$ws.Range("G89").Value = 2
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = "H"
$ws.Range("L89").Value = 1.8
$ws.Range("M89").Value = 3.6
$ws.Range("N89").Value = 4.333
$ws.Range("O89").Value = 1.75
$ws.Range("P89").Value = 3.75
$ws.Range("Q89").Value = 4.75
$ws.Range("R89").Value = -0.75
$ws.Range("S89").Value = 1.975
$ws.Range("T89").Value = 1.875
$ws.Range("U89").Value = 2.75
$ws.Range("V89").Value = 1.925
$ws.Range("W89").Value = 1.925
$ws.Range("X89").Value = 0.75
$ws.Range("Y89").Value = -1
$ws.Range("Z89").Value = -1
$ws.Range("AA89").Value = 0.9750000000000001
$ws.Range("AB89").Value = -1
$ws.Range("AC89").Value = -1
$ws.Range("AD89").Value = 0.925

# Row 170
$ws.Range("B170").Value = 6838491
$ws.Range("E170").Value = "Heracles"
$ws.Range("F170").Value = "Ajax"
$ws.Range("G170").Value = 2
$ws.Range("H170").Value = 4
$ws.Range("I170").Value = 1
$ws.Range("J170").Value = 1
$ws.Range("K170").Value = "A"
$ws.Range("L170").Value = 7.5
$ws.Range("M170").Value = 5
$ws.Range("N170").Value = 1.363
$ws.Range("O170").Value = 4.75
$ws.Range("P170").Value = 4.5
$ws.Range("Q170").Value = 1.6
$ws.Range("R170").Value = 1
$ws.Range("S170").Value = 1.85
$ws.Range("T170").Value = 2
$ws.Range("U170").Value = 3.25
$ws.Range("V170").Value = 1.875
$ws.Range("W170").Value = 1.975
$ws.Range("X170").Value = -1
$ws.Range("Y170").Value = -1
$ws.Range("Z170").Value = 0.6000000000000001
$ws.Range("AA170").Value = -1
$ws.Range("AB170").Value = 1
$ws.Range("AC170").Value = 0.875
$ws.Range("AD170").Value = -1

# Row 171
$ws.Range("B171").Value = 6838490
$ws.Range("E171").Value = "Excelsior"
$ws.Range("F171").Value = "FC Utrecht"
$ws.Range("G171").Value = 1
$ws.Range("H171").Value = 1
$ws.Range("I171").Value = 0
$ws.Range("J171").Value = 0
$ws.Range("K171").Value = "D"
$ws.Range("L171").Value = 2.7
$ws.Range("M171").Value = 3.6
$ws.Range("N171").Value = 2.375
$ws.Range("O171").Value = 3.2
$ws.Range("P171").Value = 3.6
$ws.Range("Q171").Value = 2.15
$ws.Range("R171").Value = 0.25
$ws.Range("S171").Value = 1.975
$ws.Range("T171").Value = 1.875
$ws.Range("U171").Value = 2.75
$ws.Range("V171").Value = 1.975
$ws.Range("W171").Value = 1.875
$ws.Range("X171").Value = -1
$ws.Range("Y171").Value = 2.6
$ws.Range("Z171").Value = -1
$ws.Range("AA171").Value = 0.4875
$ws.Range("AB171").Value = -0.5
$ws.Range("AC171").Value = -1
$ws.Range("AD171").Value = 0.875

# Row 279
$ws.Range("B279").Value = 7061102
$ws.Range("E279").Value = "PEC Zwolle"
$ws.Range("F279").Value = "Heracles"
$ws.Range("G279").Value = 3
$ws.Range("H279").Value = 1
$ws.Range("I279").Value = 1
$ws.Range("J279").Value = 1
$ws.Range("K279").Value = "H"
$ws.Range("L279").Value = 2.25
$ws.Range("M279").Value = 3.75
$ws.Range("N279").Value = 2.75
$ws.Range("O279").Value = 2.05
$ws.Range("P279").Value = 4
$ws.Range("Q279").Value = 3.25
$ws.Range("R279").Value = -0.5
$ws.Range("S279").Value = 2.025
$ws.Range("T279").Value = 1.825
$ws.Range("U279").Value = 3
$ws.Range("V279").Value = 1.825
$ws.Range("W279").Value = 2.025
$ws.Range("X279").Value = 1.05
$ws.Range("Y279").Value = -1
$ws.Range("Z279").Value = -1
$ws.Range("AA279").Value = 1.025
$ws.Range("AB279").Value = -1
$ws.Range("AC279").Value = 0.825
$ws.Range("AD279").Value = -1

# Row 280
$ws.Range("B280").Value = 7062784
$ws.Range("E280").Value = "NEC"
$ws.Range("F280").Value = "AZ"
$ws.Range("G280").Value = 0
$ws.Range("H280").Value = 3
$ws.Range("I280").Value = 0
$ws.Range("J280").Value = 1
$ws.Range("K280").Value = "A"
$ws.Range("L280").Value = 3.2
$ws.Range("M280").Value = 4
$ws.Range("N280").Value = 1.95
$ws.Range("O280").Value = 3.1
$ws.Range("P280").Value = 4
$ws.Range("Q280").Value = 2.1
$ws.Range("R280").Value = 0.25
$ws.Range("S280").Value = 2.03
$ws.Range("T280").Value = 1.87
$ws.Range("U280").Value = 2.75
$ws.Range("V280").Value = 1.925
$ws.Range("W280").Value = 1.925
$ws.Range("X280").Value = -1
$ws.Range("Y280").Value = -1
$ws.Range("Z280").Value = 1.1
$ws.Range("AA280").Value = -1
$ws.Range("AB280").Value = 0.8700000000000001
$ws.Range("AC280").Value = 0.4625
$ws.Range("AD280").Value = -0.5

# Row 302
$ws.Range("B302").Value = 7223358
$ws.Range("E302").Value = "Vitesse"
$ws.Range("F302").Value = "Ajax"
$ws.Range("G302").Value = 2
$ws.Range("H302").Value = 2
$ws.Range("I302").Value = 1
$ws.Range("J302").Value = 1
$ws.Range("K302").Value = "D"
$ws.Range("L302").Value = 4.75
$ws.Range("M302").Value = 4
$ws.Range("N302").Value = 1.571
$ws.Range("O302").Value = 4.333
$ws.Range("P302").Value = 4.5
$ws.Range("Q302").Value = 1.65
$ws.Range("R302").Value = 1
$ws.Range("S302").Value = 1.84
$ws.Range("T302").Value = 2.06
$ws.Range("U302").Value = 3.75
$ws.Range("V302").Value = 1.925
$ws.Range("W302").Value = 1.925
$ws.Range("X302").Value = -1
$ws.Range("Y302").Value = 3.5
$ws.Range("Z302").Value = -1
$ws.Range("AA302").Value = 0.8400000000000001
$ws.Range("AB302").Value = -1
$ws.Range("AC302").Value = 0.4625
$ws.Range("AD302").Value = -0.5

# Row 303
$ws.Range("B303").Value = 7160673
$ws.Range("E303").Value = "AZ"
$ws.Range("F303").Value = "FC Utrecht"
$ws.Range("G303").Value = 3
$ws.Range("H303").Value = 3
$ws.Range("I303").Value = 3
$ws.Range("J303").Value = 0
$ws.Range("K303").Value = "D"
$ws.Range("L303").Value = 1.4
$ws.Range("M303").Value = 4.5
$ws.Range("N303").Value = 6.5
$ws.Range("O303").Value = 1.4
$ws.Range("P303").Value = 5
$ws.Range("Q303").Value = 7
$ws.Range("R303").Value = -1.25
$ws.Range("S303").Value = 1.825
$ws.Range("T303").Value = 2.025
$ws.Range("U303").Value = 3.25
$ws.Range("V303").Value = 1.925
$ws.Range("W303").Value = 1.925
$ws.Range("X303").Value = -1
$ws.Range("Y303").Value = 4
$ws.Range("Z303").Value = -1
$ws.Range("AA303").Value = -1
$ws.Range("AB303").Value = 1.025
$ws.Range("AC303").Value = 0.925
$ws.Range("AD303").Value = -1

# Row 304
$ws.Range("B304").Value = 7155056
$ws.Range("E304").Value = "Almere City FC"
$ws.Range("F304").Value = "NEC"
$ws.Range("G304").Value = 1
$ws.Range("H304").Value = 4
$ws.Range("I304").Value = 1
$ws.Range("J304").Value = 2
$ws.Range("K304").Value = "A"
$ws.Range("L304").Value = 3.1
$ws.Range("M304").Value = 3.4
$ws.Range("N304").Value = 2.1
$ws.Range("O304").Value = 3
$ws.Range("P304").Value = 3.6
$ws.Range("Q304").Value = 2.3
$ws.Range("R304").Value = 0.25
$ws.Range("S304").Value = 1.825
$ws.Range("T304").Value = 2.025
$ws.Range("U304").Value = 2.75
$ws.Range("V304").Value = 1.8
$ws.Range("W304").Value = 2.05
$ws.Range("X304").Value = -1
$ws.Range("Y304").Value = -1
$ws.Range("Z304").Value = 1.3
$ws.Range("AA304").Value = -1
$ws.Range("AB304").Value = 1.025
$ws.Range("AC304").Value = 0.8
$ws.Range("AD304").Value = -1

# Row 305
$ws.Range("B305").Value = 7161289
$ws.Range("E305").Value = "Heracles"
$ws.Range("F305").Value = "Fortuna Sittard"
$ws.Range("G305").Value = 0
$ws.Range("H305").Value = 0
$ws.Range("I305").Value = 0
$ws.Range("J305").Value = 0
$ws.Range("K305").Value = "D"
$ws.Range("L305").Value = 2.25
$ws.Range("M305").Value = 3.5
$ws.Range("N305").Value = 2.75
$ws.Range("O305").Value = 2.6
$ws.Range("P305").Value = 3.7
$ws.Range("Q305").Value = 2.5
$ws.Range("R305").Value = 0
$ws.Range("S305").Value = 1.975
$ws.Range("T305").Value = 1.875
$ws.Range("U305").Value = 3.25
$ws.Range("V305").Value = 2.025
$ws.Range("W305").Value = 1.825
$ws.Range("X305").Value = -1
$ws.Range("Y305").Value = 2.7
$ws.Range("Z305").Value = -1
$ws.Range("AA305").Value = 0
$ws.Range("AB305").Value = 0
$ws.Range("AC305").Value = -1
$ws.Range("AD305").Value = 0.825

# Row 306
$ws.Range("B306").Value = 7223357
$ws.Range("E306").Value = "FC Volendam"
$ws.Range("F306").Value = "Go Ahead Eagles"
$ws.Range("G306").Value = 1
$ws.Range("H306").Value = 2
$ws.Range("I306").Value = 1
$ws.Range("J306").Value = 1
$ws.Range("K306").Value = "A"
$ws.Range("L306").Value = 5
$ws.Range("M306").Value = 4.333
$ws.Range("N306").Value = 1.5
$ws.Range("O306").Value = 5.25
$ws.Range("P306").Value = 4.75
$ws.Range("Q306").Value = 1.533
$ws.Range("R306").Value = 1
$ws.Range("S306").Value = 2.05
$ws.Range("T306").Value = 1.85
$ws.Range("U306").Value = 3.75
$ws.Range("V306").Value = 1.875
$ws.Range("W306").Value = 1.975
$ws.Range("X306").Value = -1
$ws.Range("Y306").Value = -1
$ws.Range("Z306").Value = 0.5329999999999999
$ws.Range("AA306").Value = 0
$ws.Range("AB306").Value = 0
$ws.Range("AC306").Value = -1
$ws.Range("AD306").Value = 0.9750000000000001
